$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fill in the missing "Taille" (height) values for the two newest players
$ws.Range("E28").Value = "1m93"
$ws.Range("E29").Value = "1m80"

# Move the active selection, matching the author's final cursor position
$ws.Range("F30").Select()
